# "hide filter in locations lists"
# Rename a few headers, append two trailing columns, resize columns,
# switch the workbook's base font to Arial, and leave the selection
# where the author last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Base font: Calibri -> Arial (the workbook's Normal / default style font)
$wb.Styles.Item("Normal").Font.Name = "Arial"

# 2. Header edits, in the order the strings were authored so the shared
#    string table lands in the same append order as the source edit.
$ws.Range("E1").Value = "Category"
$ws.Range("I1").Value = "Zone"
$ws.Range("M1").Value = "Sales Agent"
$ws.Range("L1").Value = "Service Agent"
$ws.Range("A1").Value = "Merchant ID"

# New cells should carry the same text-formatted style as the rest of
# the header row (style index 1 / numFmtId 49).
$ws.Range("L1:M1").NumberFormat = "@"

# 3. Column widths (characters, Normal-style units)
$ws.Columns.Item(1).ColumnWidth = 9.833333333333334
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 6.666666666666667
$ws.Columns.Item(4).ColumnWidth = 8.5
$ws.Columns.Item(5).ColumnWidth = 7.333333333333333
$ws.Columns.Item(6).ColumnWidth = 11.5
$ws.Columns.Item(7).ColumnWidth = 6.333333333333333
$ws.Columns.Item(8).ColumnWidth = 8
$ws.Columns.Item(9).ColumnWidth = 4
$ws.Columns.Item(10).ColumnWidth = 6.833333333333333
$ws.Columns.Item(12).ColumnWidth = 11.333333333333334
$ws.Columns.Item(13).ColumnWidth = 9.833333333333334

# 4. Selection left where the author's cursor ended up
$ws.Range("E6").Select()
